# Update column C ("Förändrad") date values from 45185 (2023-09-16) to
# 45204 (2023-10-05) for data rows 2 through 13 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45185) {
        $cell.Value = 45204
    }
}
